$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 9100.25  # H40 was 8207.429
$ws.Cells.Item(40, 10).Value = 9200  # J40 was 7612.625
$ws.Cells.Item(40, 12).Value = 9200  # L40 was 7612.625
$ws.Cells.Item(40, 14).Value = -9550  # N40 was -7962.625
$ws.Cells.Item(86, 8).Value = 5850653  # H86 was 5265887.5
$ws.Cells.Item(86, 9).Value = 3749.5  # I86 was 3599.4
$ws.Cells.Item(86, 11).Value = 3749.5  # K86 was 3599.4
$ws.Cells.Item(86, 13).Value = -2626.5  # M86 was -2476.4
$ws.Cells.Item(88, 8).Value = 4356.2  # H88 was 4597.4
$ws.Cells.Item(88, 9).Value = 0  # I88 was 5000
$ws.Cells.Item(88, 10).Value = 4356.2  # J88 was 4496.75
$ws.Cells.Item(88, 11).Value = 0  # K88 was 5000
$ws.Cells.Item(88, 12).Value = 4356.2  # L88 was 4496.75
$ws.Cells.Item(88, 13).ClearContents()  # M88 was -4594
$ws.Cells.Item(88, 14).Value = -5168.2  # N88 was -5308.75
$ws.Cells.Item(89, 8).Value = 5850653  # H89 was 5265887.5
$ws.Cells.Item(89, 9).Value = 3749.5  # I89 was 3599.4
$ws.Cells.Item(89, 11).Value = 18747.5  # K89 was 17997
$ws.Cells.Item(89, 13).Value = -13131.5  # M89 was -12381
$ws.Cells.Item(91, 8).Value = 4356.2  # H91 was 4597.4
$ws.Cells.Item(91, 9).Value = 0  # I91 was 5000
$ws.Cells.Item(91, 10).Value = 4356.2  # J91 was 4496.75
$ws.Cells.Item(91, 11).Value = 0  # K91 was 5000
$ws.Cells.Item(91, 12).Value = 4356.2  # L91 was 4496.75
$ws.Cells.Item(91, 13).ClearContents()  # M91 was -3596
$ws.Cells.Item(91, 14).Value = -7164.2  # N91 was -7304.75
$ws.Cells.Item(92, 8).Value = 160.66667  # H92 was 149.16667
$ws.Cells.Item(92, 9).Value = 160.66667  # I92 was 149.16667
$ws.Cells.Item(92, 11).Value = 160.66667  # K92 was 149.16667
$ws.Cells.Item(92, 13).Value = 1087.33333  # M92 was 1098.83333
$ws.Cells.Item(98, 8).Value = 1492.7142  # H98 was 1596.5
$ws.Cells.Item(98, 9).Value = 1190.119  # I98 was 1302.2559
$ws.Cells.Item(98, 10).Value = 2400.5  # J98 was 2569.7693
$ws.Cells.Item(98, 11).Value = 1190.119  # K98 was 1302.2559
$ws.Cells.Item(98, 12).Value = 2400.5  # L98 was 2569.7693
$ws.Cells.Item(98, 13).Value = 307.8810000000001  # M98 was 195.7440999999999
$ws.Cells.Item(98, 14).Value = -5396.5  # N98 was -5565.7693
$ws.Cells.Item(111, 8).Value = 61531.53  # H111 was 65157.312
$ws.Cells.Item(111, 9).Value = 102438.5  # I111 was 113429.555
$ws.Cells.Item(111, 11).Value = 307315.5  # K111 was 340288.665
$ws.Cells.Item(111, 13).Value = -304248.5  # M111 was -337221.665
$ws.Cells.Item(122, 8).Value = 1492.7142  # H122 was 1596.5
$ws.Cells.Item(122, 9).Value = 1190.119  # I122 was 1302.2559
$ws.Cells.Item(122, 10).Value = 2400.5  # J122 was 2569.7693
$ws.Cells.Item(122, 11).Value = 3570.357  # K122 was 3906.7677
$ws.Cells.Item(122, 12).Value = 7201.5  # L122 was 7709.3079
$ws.Cells.Item(122, 13).Value = -1120.357  # M122 was -1456.7677
$ws.Cells.Item(122, 14).Value = -12101.5  # N122 was -12609.3079
$ws.Cells.Item(132, 8).Value = 1843.4225  # H132 was 1984.8209
$ws.Cells.Item(132, 9).Value = 1930.0769  # I132 was 2091.0657
$ws.Cells.Item(132, 11).Value = 5790.2307  # K132 was 6273.1971
$ws.Cells.Item(132, 13).Value = -3260.2307  # M132 was -3743.1971
$ws.Cells.Item(138, 8).Value = 6152.0127  # H138 was 6174.859
$ws.Cells.Item(138, 10).Value = 6801.183  # J138 was 6854.293
$ws.Cells.Item(138, 12).Value = 20403.549  # L138 was 20562.879
$ws.Cells.Item(138, 14).Value = -30683.549  # N138 was -30842.879
$ws.Cells.Item(141, 8).Value = 4132.28  # H141 was 2949.5
$ws.Cells.Item(141, 9).Value = 1925.5  # I141 was 1376.2693
$ws.Cells.Item(141, 10).Value = 8055.4443  # J141 was 8062.5
$ws.Cells.Item(141, 11).Value = 5776.5  # K141 was 4128.8079
$ws.Cells.Item(141, 12).Value = 24166.3329  # L141 was 24187.5
$ws.Cells.Item(141, 13).Value = -596.5  # M141 was 1051.1921
$ws.Cells.Item(141, 14).Value = -34526.3329  # N141 was -34547.5

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1097.98  # H32 was 1189.2
$ws.Cells.Item(32, 9).Value = 1097.98  # I32 was 1189.2
$ws.Cells.Item(32, 11).Value = 1097.98  # K32 was 1189.2
$ws.Cells.Item(32, 13).Value = -810.98  # M32 was -902.2
$ws.Cells.Item(61, 8).Value = 2709.322  # H61 was 2695.2307
$ws.Cells.Item(61, 9).Value = 1855.0408  # I61 was 1817.88
$ws.Cells.Item(61, 10).Value = 6895.3  # J61 was 5619.7334
$ws.Cells.Item(61, 11).Value = 1855.0408  # K61 was 1817.88
$ws.Cells.Item(61, 12).Value = 6895.3  # L61 was 5619.7334
$ws.Cells.Item(61, 13).Value = -1643.0408  # M61 was -1605.88
$ws.Cells.Item(61, 14).Value = -7319.3  # N61 was -6043.7334
$ws.Cells.Item(74, 8).Value = 188333.33  # H74 was 0
$ws.Cells.Item(74, 9).Value = 132500  # I74 was 0
$ws.Cells.Item(74, 10).Value = 300000  # J74 was 0
$ws.Cells.Item(74, 11).Value = 132500  # K74 was 0
$ws.Cells.Item(74, 12).Value = 300000  # L74 was 0
$ws.Cells.Item(74, 13).Value = -131626  # M74 was None
$ws.Cells.Item(74, 14).Value = -301748  # N74 was None
$ws.Cells.Item(77, 8).Value = 188333.33  # H77 was 0
$ws.Cells.Item(77, 9).Value = 132500  # I77 was 0
$ws.Cells.Item(77, 10).Value = 300000  # J77 was 0
$ws.Cells.Item(77, 11).Value = 662500  # K77 was 0
$ws.Cells.Item(77, 12).Value = 1500000  # L77 was 0
$ws.Cells.Item(77, 13).Value = -658132  # M77 was None
$ws.Cells.Item(77, 14).Value = -1508736  # N77 was None
$ws.Cells.Item(95, 8).Value = 34104  # H95 was 108534.336
$ws.Cells.Item(95, 10).Value = 34104  # J95 was 108534.336
$ws.Cells.Item(95, 12).Value = 34104  # L95 was 108534.336
$ws.Cells.Item(95, 14).Value = -39596  # N95 was -114026.336
$ws.Cells.Item(103, 8).Value = 73000  # H103 was 64249.668
$ws.Cells.Item(103, 10).Value = 73000  # J103 was 64249.668
$ws.Cells.Item(103, 12).Value = 73000  # L103 was 64249.668
$ws.Cells.Item(103, 14).Value = -75344  # N103 was -66593.66800000001
$ws.Cells.Item(132, 8).Value = 3889.5305  # H132 was 4105.7393
$ws.Cells.Item(132, 9).Value = 2661.9666  # I132 was 2951.423
$ws.Cells.Item(132, 10).Value = 5827.7896  # J132 was 5606.35
$ws.Cells.Item(132, 11).Value = 7985.899800000001  # K132 was 8854.269
$ws.Cells.Item(132, 12).Value = 17483.3688  # L132 was 16819.05
$ws.Cells.Item(132, 13).Value = -5455.899800000001  # M132 was -6324.269
$ws.Cells.Item(132, 14).Value = -22543.3688  # N132 was -21879.05
$ws.Cells.Item(136, 8).Value = 2709.322  # H136 was 2695.2307
$ws.Cells.Item(136, 9).Value = 1855.0408  # I136 was 1817.88
$ws.Cells.Item(136, 10).Value = 6895.3  # J136 was 5619.7334
$ws.Cells.Item(136, 11).Value = 5565.1224  # K136 was 5453.64
$ws.Cells.Item(136, 12).Value = 20685.9  # L136 was 16859.2002
$ws.Cells.Item(136, 13).Value = -3015.1224  # M136 was -2903.64
$ws.Cells.Item(136, 14).Value = -25785.9  # N136 was -21959.2002

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 42138.25  # H76 was 45184.668
$ws.Cells.Item(76, 10).Value = 42138.25  # J76 was 45184.668
$ws.Cells.Item(76, 12).Value = 42138.25  # L76 was 45184.668
$ws.Cells.Item(76, 14).Value = -42768.25  # N76 was -45814.668
$ws.Cells.Item(79, 8).Value = 42138.25  # H79 was 45184.668
$ws.Cells.Item(79, 10).Value = 42138.25  # J79 was 45184.668
$ws.Cells.Item(79, 12).Value = 42138.25  # L79 was 45184.668
$ws.Cells.Item(79, 14).Value = -44322.25  # N79 was -47368.668
$ws.Cells.Item(99, 8).Value = 5293.1  # H99 was 5398.3687
$ws.Cells.Item(99, 9).Value = 5413.25  # I99 was 5554.6
$ws.Cells.Item(99, 11).Value = 5413.25  # K99 was 5554.6
$ws.Cells.Item(99, 13).Value = -3915.25  # M99 was -4056.6
$ws.Cells.Item(105, 8).Value = 1332.3846  # H105 was 1301.7142
$ws.Cells.Item(105, 9).Value = 1110.125  # I105 was 1094.1923
$ws.Cells.Item(105, 11).Value = 1110.125  # K105 was 1094.1923
$ws.Cells.Item(105, 13).Value = 636.875  # M105 was 652.8077000000001
$ws.Cells.Item(134, 8).Value = 15529.278  # H134 was 15382.475
$ws.Cells.Item(134, 9).Value = 1865.3729  # I134 was 1897.2931
$ws.Cells.Item(134, 10).Value = 55837.8  # J134 was 50934.316
$ws.Cells.Item(134, 11).Value = 5596.1187  # K134 was 5691.879300000001
$ws.Cells.Item(134, 12).Value = 167513.4  # L134 was 152802.948
$ws.Cells.Item(134, 13).Value = -3061.1187  # M134 was -3156.879300000001
$ws.Cells.Item(134, 14).Value = -172583.4  # N134 was -157872.948

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 192689  # H58 was 192675.7
$ws.Cells.Item(58, 9).Value = 387252.44  # I58 was 372920.8
$ws.Cells.Item(58, 10).Value = 5331.6294  # J58 was 5498.077
$ws.Cells.Item(58, 11).Value = 387252.44  # K58 was 372920.8
$ws.Cells.Item(58, 12).Value = 5331.6294  # L58 was 5498.077
$ws.Cells.Item(58, 13).Value = -387049.44  # M58 was -372717.8
$ws.Cells.Item(58, 14).Value = -5737.6294  # N58 was -5904.077
$ws.Cells.Item(122, 8).Value = 3427.96  # H122 was 3801.5908
$ws.Cells.Item(122, 9).Value = 2696.077  # I122 was 3106
$ws.Cells.Item(122, 10).Value = 4220.8335  # J122 was 4381.25
$ws.Cells.Item(122, 11).Value = 8088.231000000001  # K122 was 9318
$ws.Cells.Item(122, 12).Value = 12662.5005  # L122 was 13143.75
$ws.Cells.Item(122, 13).Value = -5638.231000000001  # M122 was -6868
$ws.Cells.Item(122, 14).Value = -17562.5005  # N122 was -18043.75
$ws.Cells.Item(132, 8).Value = 3204.4285  # H132 was 2435.93
$ws.Cells.Item(132, 9).Value = 2073.9062  # I132 was 1594.2709
$ws.Cells.Item(132, 10).Value = 6822.1  # J132 was 6924.778
$ws.Cells.Item(132, 11).Value = 6221.7186  # K132 was 4782.8127
$ws.Cells.Item(132, 12).Value = 20466.3  # L132 was 20774.334
$ws.Cells.Item(132, 13).Value = -3691.7186  # M132 was -2252.8127
$ws.Cells.Item(132, 14).Value = -25526.3  # N132 was -25834.334
$ws.Cells.Item(134, 8).Value = 199720.16  # H134 was 199732.39
$ws.Cells.Item(134, 9).Value = 2570.8386  # I134 was 2721.7585
$ws.Cells.Item(134, 10).Value = 505301.6  # J134 was 459428.22
$ws.Cells.Item(134, 11).Value = 7712.5158  # K134 was 8165.2755
$ws.Cells.Item(134, 12).Value = 1515904.8  # L134 was 1378284.66
$ws.Cells.Item(134, 13).Value = -5177.5158  # M134 was -5630.2755
$ws.Cells.Item(134, 14).Value = -1520974.8  # N134 was -1383354.66
$ws.Cells.Item(136, 8).Value = 192689  # H136 was 192675.7
$ws.Cells.Item(136, 9).Value = 387252.44  # I136 was 372920.8
$ws.Cells.Item(136, 10).Value = 5331.6294  # J136 was 5498.077
$ws.Cells.Item(136, 11).Value = 1161757.32  # K136 was 1118762.4
$ws.Cells.Item(136, 12).Value = 15994.8882  # L136 was 16494.231
$ws.Cells.Item(136, 13).Value = -1159207.32  # M136 was -1116212.4
$ws.Cells.Item(136, 14).Value = -21094.8882  # N136 was -21594.231

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 45073424  # H4 was 41469564
$ws.Cells.Item(4, 9).Value = 57537572  # I4 was 57537580
$ws.Cells.Item(4, 10).Value = 202494.4  # J4 was 151806.58
$ws.Cells.Item(4, 11).Value = 172612716  # K4 was 172612740
$ws.Cells.Item(4, 12).Value = 607483.2  # L4 was 455419.74
$ws.Cells.Item(4, 13).Value = -172612604  # M4 was -172612628
$ws.Cells.Item(4, 14).Value = -607707.2  # N4 was -455643.74
$ws.Cells.Item(14, 8).Value = 392.14285  # H14 was 352.125
$ws.Cells.Item(14, 9).Value = 392.14285  # I14 was 352.125
$ws.Cells.Item(14, 11).Value = 1176.42855  # K14 was 1056.375
$ws.Cells.Item(14, 13).Value = -1003.42855  # M14 was -883.375
$ws.Cells.Item(17, 8).Value = 833.3333  # H17 was 759.6
$ws.Cells.Item(17, 9).Value = 925  # I17 was 833
$ws.Cells.Item(17, 10).Value = 650  # J17 was 649.5
$ws.Cells.Item(17, 11).Value = 2775  # K17 was 2499
$ws.Cells.Item(17, 12).Value = 1950  # L17 was 1948.5
$ws.Cells.Item(17, 13).Value = -2606  # M17 was -2330
$ws.Cells.Item(17, 14).Value = -2288  # N17 was -2286.5
$ws.Cells.Item(129, 8).Value = 1765.8  # H129 was 1999.2354
$ws.Cells.Item(129, 9).Value = 519.3333  # I129 was 534.25
$ws.Cells.Item(129, 10).Value = 3635.5  # J129 was 3301.4443
$ws.Cells.Item(129, 11).Value = 1557.9999  # K129 was 1602.75
$ws.Cells.Item(129, 12).Value = 10906.5  # L129 was 9904.332900000001
$ws.Cells.Item(129, 13).Value = 3442.0001  # M129 was 3397.25
$ws.Cells.Item(129, 14).Value = -20906.5  # N129 was -19904.3329

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 7499999.5  # H7 was 4126250
$ws.Cells.Item(7, 9).Value = 35000000  # I7 was 0
$ws.Cells.Item(7, 10).Value = 2916666  # J7 was 4126250
$ws.Cells.Item(7, 11).Value = 35000000  # K7 was 0
$ws.Cells.Item(7, 12).Value = 2916666  # L7 was 4126250
$ws.Cells.Item(7, 13).Value = -34999888  # M7 was None
$ws.Cells.Item(7, 14).Value = -2916890  # N7 was -4126474
$ws.Cells.Item(8, 8).Value = 7499999.5  # H8 was 4126250
$ws.Cells.Item(8, 9).Value = 35000000  # I8 was 0
$ws.Cells.Item(8, 10).Value = 2916666  # J8 was 4126250
$ws.Cells.Item(8, 11).Value = 35000000  # K8 was 0
$ws.Cells.Item(8, 12).Value = 2916666  # L8 was 4126250
$ws.Cells.Item(8, 13).Value = -34999861  # M8 was None
$ws.Cells.Item(8, 14).Value = -2916944  # N8 was -4126528
$ws.Cells.Item(109, 8).Value = 58999  # H109 was 56999.5
$ws.Cells.Item(109, 10).Value = 58999  # J109 was 56999.5
$ws.Cells.Item(109, 12).Value = 58999  # L109 was 56999.5
$ws.Cells.Item(109, 14).Value = -61079  # N109 was -59079.5
$ws.Cells.Item(122, 8).Value = 2714.125  # H122 was 2925.4333
$ws.Cells.Item(122, 9).Value = 2267.7778  # I122 was 2608.1875
$ws.Cells.Item(122, 11).Value = 6803.3334  # K122 was 7824.5625
$ws.Cells.Item(122, 13).Value = -4353.3334  # M122 was -5374.5625
$ws.Cells.Item(141, 8).Value = 28816.334  # H141 was 29579.6
$ws.Cells.Item(141, 10).Value = 28816.334  # J141 was 29579.6
$ws.Cells.Item(141, 12).Value = 28816.334  # L141 was 29579.6
$ws.Cells.Item(141, 14).Value = -39176.334  # N141 was -39939.6

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 1558001.2  # H20 was 1557001.2
$ws.Cells.Item(20, 9).Value = 3750000  # I20 was 2505000
$ws.Cells.Item(20, 10).Value = 96668.664  # J20 was 135003
$ws.Cells.Item(20, 11).Value = 3750000  # K20 was 2505000
$ws.Cells.Item(20, 12).Value = 96668.664  # L20 was 135003
$ws.Cells.Item(20, 13).Value = -3749774  # M20 was -2504774
$ws.Cells.Item(20, 14).Value = -97120.664  # N20 was -135455
$ws.Cells.Item(39, 8).Value = 15000  # H39 was 0
$ws.Cells.Item(39, 10).Value = 15000  # J39 was 0
$ws.Cells.Item(39, 12).Value = 15000  # L39 was 0
$ws.Cells.Item(39, 13).Value = 0  # M39 was None
$ws.Cells.Item(39, 14).Value = -15920  # N39 was None
$ws.Cells.Item(46, 8).Value = 5272.636  # H46 was 5499.8667
$ws.Cells.Item(46, 9).Value = 5000  # I46 was 4999.9
$ws.Cells.Item(46, 10).Value = 6499.5  # J46 was 6499.8
$ws.Cells.Item(46, 11).Value = 5000  # K46 was 4999.9
$ws.Cells.Item(46, 12).Value = 6499.5  # L46 was 6499.8
$ws.Cells.Item(46, 13).Value = -4812  # M46 was -4811.9
$ws.Cells.Item(46, 14).Value = -6875.5  # N46 was -6875.8
$ws.Cells.Item(122, 8).Value = 957777  # H122 was 1182171
$ws.Cells.Item(122, 9).Value = 914119.4  # I122 was 1255362.9
$ws.Cells.Item(122, 10).Value = 1005800.4  # J122 was 1117111.6
$ws.Cells.Item(122, 11).Value = 2742358.2  # K122 was 3766088.7
$ws.Cells.Item(122, 12).Value = 3017401.2  # L122 was 3351334.8
$ws.Cells.Item(122, 13).Value = -2739908.2  # M122 was -3763638.7
$ws.Cells.Item(122, 14).Value = -3022301.2  # N122 was -3356234.8
$ws.Cells.Item(123, 8).Value = 76982.25  # H123 was 76989.25
$ws.Cells.Item(123, 10).Value = 76982.25  # J123 was 76989.25
$ws.Cells.Item(123, 12).Value = 76982.25  # L123 was 76989.25
$ws.Cells.Item(123, 14).Value = -86782.25  # N123 was -86789.25

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 0  # H54 was 9000
$ws.Cells.Item(54, 9).Value = 0  # I54 was 9000
$ws.Cells.Item(54, 11).Value = 0  # K54 was 9000
$ws.Cells.Item(54, 13).ClearContents()  # M54 was -8480
$ws.Cells.Item(100, 8).Value = 1175.25  # H100 was 1250
$ws.Cells.Item(100, 9).Value = 1100.5  # I100 was 0
$ws.Cells.Item(100, 11).Value = 2201  # K100 was 0
$ws.Cells.Item(100, 13).Value = -1660  # M100 was None
$ws.Cells.Item(110, 8).Value = 46997.5  # H110 was 48497
$ws.Cells.Item(110, 10).Value = 46997.5  # J110 was 48497
$ws.Cells.Item(110, 12).Value = 46997.5  # L110 was 48497
$ws.Cells.Item(110, 14).Value = -55177.5  # N110 was -56677
$ws.Cells.Item(122, 8).Value = 30306104  # H122 was 28574332
$ws.Cells.Item(122, 9).Value = 45456224  # I122 was 41668216
$ws.Cells.Item(122, 11).Value = 136368672  # K122 was 125004648
$ws.Cells.Item(122, 13).Value = -136366222  # M122 was -125002198
$ws.Cells.Item(132, 8).Value = 17717.016  # H132 was 18013.5
$ws.Cells.Item(132, 9).Value = 2196.173  # I132 was 2250.2307
$ws.Cells.Item(132, 10).Value = 71522.60000000001  # J132 was 76562.78999999999
$ws.Cells.Item(132, 11).Value = 6588.518999999999  # K132 was 6750.6921
$ws.Cells.Item(132, 12).Value = 214567.8  # L132 was 229688.37
$ws.Cells.Item(132, 13).Value = -4058.518999999999  # M132 was -4220.6921
$ws.Cells.Item(132, 14).Value = -219627.8  # N132 was -234748.37
